$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 60 (before the current row 61),
# pushing the existing rows 61-67 down to 63-69.
$ws.Rows.Item(61).Insert()
$ws.Rows.Item(61).Insert()

# New row 61: week of 2023-04-18 (serial 45034), quality "Primera"
$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 45034
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100104
$ws.Range("H61").Value = "Frutos de pepita"
$ws.Range("I61").Value = 100104003
$ws.Range("J61").Value = "Membrillo"
$ws.Range("K61").Value = "Champion"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 16
$ws.Range("N61").Value = 270000
$ws.Range("O61").Value = 280000
$ws.Range("P61").Value = 275000
$ws.Range("Q61").Value = "`$/bins (450 kilos)"
$ws.Range("R61").Value = "Región de O'Higgins"
$ws.Range("S61").Value = 611
$ws.Range("T61").Value = 450

# New row 62: week of 2023-04-18 (serial 45034), quality "Segunda"
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 45034
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100104
$ws.Range("H62").Value = "Frutos de pepita"
$ws.Range("I62").Value = 100104003
$ws.Range("J62").Value = "Membrillo"
$ws.Range("K62").Value = "Champion"
$ws.Range("L62").Value = "Segunda"
$ws.Range("M62").Value = 10
$ws.Range("N62").Value = 230000
$ws.Range("O62").Value = 240000
$ws.Range("P62").Value = 235000
$ws.Range("Q62").Value = "`$/bins (450 kilos)"
$ws.Range("R62").Value = "Región de O'Higgins"
$ws.Range("S62").Value = 522
$ws.Range("T62").Value = 450

# Make sure the date cells keep the same date-time number format as the
# rest of column D (style index 2 in the original workbook).
$ws.Range("D61").NumberFormat = $ws.Range("D60").NumberFormat
$ws.Range("D62").NumberFormat = $ws.Range("D60").NumberFormat
